# Applies the crypto price/volume refresh described in the commit diff.
# Values in columns B-E are plain text (coinranking.com export format),
# so each cell is forced to Text format before the write and restored to
# the default ('Normal') style afterwards -- this stops Excel's COM layer
# from auto-coercing number-looking strings (e.g. '0.9995', '30.520.27')
# into floating point values, which would corrupt the display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '30.520.27'
Set-TextValue 'E2' '  -1.42%  '
Set-TextValue 'D3' '1.912.97'
Set-TextValue 'E3' '  -2.13%  '
Set-TextValue 'D4' '0.9995'
Set-TextValue 'E4' '  -0.24%  '
Set-TextValue 'D5' '239.69'
Set-TextValue 'E5' '  -1.42%  '
Set-TextValue 'D6' '0.9997'
Set-TextValue 'E6' '  -0.18%  '
Set-TextValue 'D7' '0.4753'
Set-TextValue 'E7' '  -2.22%  '
Set-TextValue 'D8' '0.2844'
Set-TextValue 'E8' '  -3.00%  '
Set-TextValue 'D9' '0.06691'
Set-TextValue 'E9' '  -4.51%  '
Set-TextValue 'D10' '18.80'
Set-TextValue 'E10' '  -3.33%  '
Set-TextValue 'D11' '101.23'
Set-TextValue 'E11' '  -5.80%  '
Set-TextValue 'D12' '1.921.86'
Set-TextValue 'E12' '  -2.28%  '
Set-TextValue 'D13' '0.07682'
Set-TextValue 'D14' '5.235'
Set-TextValue 'E14' '  -2.15%  '
Set-TextValue 'D15' '0.6703'
Set-TextValue 'E15' '  -4.00%  '
Set-TextValue 'D16' '30.533.61'
Set-TextValue 'E16' '  -1.44%  '
Set-TextValue 'D17' '256.25'
Set-TextValue 'E17' '  -7.67%  '
Set-TextValue 'D18' '0.9984'
Set-TextValue 'E18' '  -0.28%  '
Set-TextValue 'D19' '0.000007485'
Set-TextValue 'E19' '  -3.35%  '
Set-TextValue 'E20' '  -3.97%  '
Set-TextValue 'D21' '5.404'
Set-TextValue 'E21' '  -1.17%  '
Set-TextValue 'D22' '0.9996'
Set-TextValue 'E22' '  -0.31%  '
Set-TextValue 'D23' '0.4522'
Set-TextValue 'E23' '  -9.41%  '
Set-TextValue 'D24' '6.304'
Set-TextValue 'E24' '  -2.76%  '
Set-TextValue 'D25' '168.51'
Set-TextValue 'E25' '  -0.02%  '
Set-TextValue 'D26' '9.358'
Set-TextValue 'E26' '  -3.93%  '
Set-TextValue 'D27' '18.99'
Set-TextValue 'E27' '  -3.30%  '
Set-TextValue 'E28' '  -4.83%  '
Set-TextValue 'D29' '4.717'
Set-TextValue 'E29' '  +2.96%  '
Set-TextValue 'D30' '0.1007'
Set-TextValue 'E30' '  -3.56%  '
Set-TextValue 'D31' '1.380'
Set-TextValue 'E31' '  -1.66%  '
Set-TextValue 'E32' '  -3.09%  '
Set-TextValue 'D33' '4.261'
Set-TextValue 'E33' '  -3.00%  '
Set-TextValue 'D34' '0.04730'
Set-TextValue 'E34' '  -2.92%  '
Set-TextValue 'D35' '0.7286'
Set-TextValue 'E35' '  -2.90%  '
Set-TextValue 'D36' '1.112'
Set-TextValue 'E36' '  -4.30%  '
Set-TextValue 'D37' '0.9981'
Set-TextValue 'E37' '  -0.22%  '
Set-TextValue 'D38' '2.718'
Set-TextValue 'E38' '  -0.51%  '
Set-TextValue 'D39' '0.01917'
Set-TextValue 'E39' '  -3.79%  '
Set-TextValue 'D40' '2.615'
Set-TextValue 'E40' '  -2.35%  '
Set-TextValue 'B41' 'FraxShare'
Set-TextValue 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D41' '6.232'
Set-TextValue 'E41' '  -4.40%  '
Set-TextValue 'B42' 'Aave'
Set-TextValue 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '74.85'
Set-TextValue 'E42' '  -4.25%  '
Set-TextValue 'D43' '1.956'
Set-TextValue 'E43' '  -7.06%  '
Set-TextValue 'D44' '0.8614'
Set-TextValue 'E44' '  -3.77%  '
Set-TextValue 'D45' '105.30'
Set-TextValue 'E45' '  -3.53%  '
Set-TextValue 'D46' '0.4243'
Set-TextValue 'E46' '  -4.16%  '
Set-TextValue 'D47' '0.9992'
Set-TextValue 'E47' '  -0.12%  '
Set-TextValue 'D48' '988.67'
Set-TextValue 'E48' '  +0.06%  '
Set-TextValue 'D49' '7.410'
Set-TextValue 'E49' '  -4.81%  '
Set-TextValue 'D50' '0.1199'
Set-TextValue 'E50' '  -3.94%  '
Set-TextValue 'D51' '34.80'
Set-TextValue 'E51' '  -3.02%  '
